$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text format so numeric-looking strings are not
# auto-converted to numbers (matches existing inlineStr text cells).
$valueRange = $ws.Range("D2:E51")
$valueRange.NumberFormat = "@"

# Row 25 and 26 swap: Monero <-> Toncoin (full row content swap)
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "1.851"
$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "149.57"
$ws.Range("E26").Value = "  -2.26%  "

# Price / Volume updates
$ws.Range("D2").Value = "27.088.04"
$ws.Range("D3").Value = "1.867.54"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "307.02"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.5087"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("D8").Value = "0.3742"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "0.07134"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("D10").Value = "0.8849"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").Value = "20.60"
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("D12").Value = "1.900.66"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "0.07525"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "89.06"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").Value = "0.9999"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "0.000008463"
$ws.Range("E17").Value = "  -3.12%  "
$ws.Range("D18").Value = "14.10"
$ws.Range("E18").Value = "  -3.74%  "
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "27.138.54"
$ws.Range("E20").Value = "  -2.78%  "
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").Value = "2.113.52"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("D24").Value = "6.471"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").Value = "2.098"
$ws.Range("E28").Value = "  -4.81%  "
$ws.Range("D29").Value = "112.77"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Value = "4.732"
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("D31").Value = "4.682"
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("D32").Value = "0.09004"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "0.05124"
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("D34").Value = "3.088"
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("E35").Value = "  -6.03%  "
$ws.Range("D36").Value = "0.7329"
$ws.Range("E36").Value = "  -5.96%  "
$ws.Range("D37").Value = "0.02046"
$ws.Range("E37").Value = "  -1.72%  "
$ws.Range("D38").Value = "2.502"
$ws.Range("E38").Value = "  -4.74%  "
$ws.Range("D39").Value = "3.059"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").Value = "0.5316"
$ws.Range("E41").Value = "  -4.16%  "
$ws.Range("D42").Value = "6.558"
$ws.Range("E42").Value = "  -3.90%  "
$ws.Range("D43").Value = "117.08"
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("D44").Value = "8.309"
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("D46").Value = "0.9996"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "0.4617"
$ws.Range("E47").Value = "  -4.12%  "
$ws.Range("D48").Value = "9.980"
$ws.Range("E48").Value = "  -6.03%  "
$ws.Range("D49").Value = "1.561"
$ws.Range("E49").Value = "  -4.80%  "
$ws.Range("D50").Value = "64.38"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("D51").Value = "36.45"
$ws.Range("E51").Value = "  -1.67%  "

# Restore the default (Normal) style so no stray NumberFormat/style
# metadata is left behind on these cells.
$valueRange.Style = "Normal"
